# Update countries & provincias Spain
# Applies updated case-data values for several countries. One of the
# updates (Sudan) increases its "Casos totales" (column B) enough that it
# now ranks above Vietnam and Guinea Ecuatorial in the (descending sorted)
# list, so those two rows' data shift down by one row and Sudan's row gets
# fresh figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 58: Moldavia -------------------------------------------------
$ws.Range("B58").Value2 = 3481
$ws.Range("C58").Value2 = 73
$ws.Range("E58").Value2 = 2455

# --- Row 120: Mauricio -------------------------------------------------
$ws.Range("B120").Value2 = 334
$ws.Range("C120").Value2 = 2
$ws.Range("D120").Value2 = 302
$ws.Range("E120").Value2 = 23

# --- Rows 126-128: Sudan moves above Vietnam / Guinea Ecuatorial -------
# Row 126 now holds Sudan's refreshed figures.
$ws.Range("A126").Value2 = "Sudan"
$ws.Range("B126").Value2 = 275
$ws.Range("C126").Value2 = 38
$ws.Range("D126").Value2 = 21
$ws.Range("E126").Value2 = 232
$ws.Range("F126").Value2 = 0
$ws.Range("G126").Value2 = 1
$ws.Range("H126").Value2 = 22

# Row 127 now holds what used to be Vietnam's row (shifted down).
$ws.Range("A127").Value2 = "Vietnam"
$ws.Range("B127").Value2 = 270
$ws.Range("C127").Value2 = 0
$ws.Range("D127").Value2 = 225
$ws.Range("E127").Value2 = 45
$ws.Range("F127").Value2 = 8
$ws.Range("G127").Value2 = 0
$ws.Range("H127").Value2 = 0

# Row 128 now holds what used to be Guinea Ecuatorial's row (shifted down).
$ws.Range("A128").Value2 = "Guinea Ecuatorial"
$ws.Range("B128").Value2 = 258
$ws.Range("C128").Value2 = 0
$ws.Range("D128").Value2 = 8
$ws.Range("E128").Value2 = 249
$ws.Range("F128").Value2 = 0
$ws.Range("G128").Value2 = 0
$ws.Range("H128").Value2 = 1

# --- Row 144: Trinidad y Tobago ----------------------------------------
$ws.Range("D144").Value2 = 59
$ws.Range("E144").Value2 = 49
